# Quarterly indexing esoteric bug-fix operation
#
# Column A (rows 2-73) holds quarter-start date serials (the 1st of
# Jan/Apr/Jul/Oct). The fix re-indexes each one to the 15th of the
# *following* month (Feb/May/Aug/Nov 15th), i.e.:
#     new = DATE(YEAR(old), MONTH(old) + 1, 15)
#
# Note: `$cell.Value` on a date-formatted cell returns an opaque COM
# Variant wrapper in this host that won't stringify/coerce cleanly, so
# we read/write the raw numeric date serial via `.Value2` instead, and
# do the calendar math through .NET DateTime (`Get-Date` / `.AddDays`)
# converted back to an OLE Automation date serial via `.ToOADate()`.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$epoch = Get-Date -Year 1899 -Month 12 -Day 30

for ($r = 2; $r -le 73; $r++) {
    $cell = $ws.Cells.Item($r, 1)

    $serial = $cell.Value2
    $oldDate = $epoch.AddDays($serial)

    $nextMonthDate = $oldDate.AddMonths(1)
    $newDate = Get-Date -Year $nextMonthDate.Year -Month $nextMonthDate.Month -Day 15

    $cell.Value2 = [Math]::Floor($newDate.ToOADate())
}
